$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 371; this shifts the existing rows 371-395 down to 372-396,
# preserving all of their data/styles exactly, and leaves row 371 blank
# (inheriting the date-column style from the row above, like Excel normally does).
$ws.Rows.Item(371).Insert()

# Populate the new record in row 371.
$ws.Cells.Item(371, 1).Value = 9
$ws.Cells.Item(371, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(371, 3).Value = "Metropolitana"
$ws.Cells.Item(371, 4).Value = (Get-Date -Year 2023 -Month 4 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(371, 5).Value = 13
$ws.Cells.Item(371, 6).Value = 100112030
$ws.Cells.Item(371, 7).Value = "Poroto granado"
$ws.Cells.Item(371, 8).Value = "Sin especificar"
$ws.Cells.Item(371, 9).Value = "Primera"
$ws.Cells.Item(371, 10).Value = 70
$ws.Cells.Item(371, 11).Value = 31000
$ws.Cells.Item(371, 12).Value = 34000
$ws.Cells.Item(371, 13).Value = 32500
$ws.Cells.Item(371, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(371, 15).Value = "Región Metropolitana"
$ws.Cells.Item(371, 16).Value = 1300
$ws.Cells.Item(371, 17).Value = 25
$ws.Cells.Item(371, 18).Value = "Hortaliza"
